$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "ddd"
$ws.Range("A5").Value = "ddd"

$ws.Range("A6").Select()
